# Refactor synthetic array: insert a new "statut_name" column at column C,
# shifting the existing NCTId..intervention_type columns one position right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C (NCTId), pushing
# everything from C onward (now D onward) to the right.
$ws.Columns.Item(3).Insert()

# New header in C1, matching style of the other header cells (copy from B1).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C1").Value = "statut_name"

# Fill in the statut_name values for each data row (2-16).
$statutNames = @{
    2  = "pas de résultat ni de publication"
    3  = "pas de résultat ni de publication"
    4  = "résultat et / ou publication posté"
    5  = "résultat et / ou publication posté"
    6  = "pas de résultat ni de publication"
    7  = "pas de résultat ni de publication"
    8  = "résultat et / ou publication posté dans les 12 mois"
    9  = "pas de résultat ni de publication"
    10 = "pas de résultat ni de publication"
    11 = "pas de résultat ni de publication"
    12 = "pas de résultat ni de publication"
    13 = "pas de résultat ni de publication"
    14 = "pas de résultat ni de publication"
    15 = "pas de résultat ni de publication"
    16 = "pas de résultat ni de publication"
}

foreach ($row in $statutNames.Keys) {
    $ws.Cells.Item($row, 3).Value = $statutNames[$row]
}

$excel.CutCopyMode = $false
